# Goal: document a new TODO note in the "Goal" worksheet (rows 20-21),
# matching the commit "add feature (save Item table) only one row. need to
# develope adding row."
#
# The new note is inserted right after the existing 2017-03-31 entry (row 18)
# and before the "장기 목표" (long-term goal) block, which therefore shifts
# down by the same number of rows that are inserted.

$wb = $excel.ActiveWorkbook

$goal = $wb.Worksheets.Item(8)     # "Goal"
$case = $wb.Worksheets.Item(7)     # "Case Table"

# --- Case Table: last selection moves to B8 (recorded while not the active sheet) ---
$case.Activate()
$case.Range("B8").Select()

# --- Goal sheet: insert 9 blank rows at row 19 so the old row 24 block
# ("장기 목표 " and everything below it) ends up at row 33, matching the
# new dimension A3:H42 ---
$goal.Activate()
$goal.Rows("19:27").Insert()

# New short-term TODO entry (row 20/21)
$goal.Range("A20").Value = 42846
$goal.Range("B20").Value = "1."
$goal.Range("C20").Value = "company와 site쪽에도 검색 항목 바로 완성이아니라 고를 수있도록 구현.. "
$goal.Range("B21").Value = "2."
$goal.Range("C21").Value = "한글 저장안됨.."

$goal.Rows("20:20").RowHeight = 34.8

# Restore the view/selection state recorded for the Goal sheet
$goal.Range("C21").Select()
